$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right.
$ws.Columns.Item(1).Insert()

# ---- Write new-content cells first, in the exact order that reproduces the
# ---- shared-string table ordering of the target workbook. ----

$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "FilesTab"

$countQuery1 = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE f.file_type = 'Aligned DNA reads file'
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@
$ws.Range("C2").Value = $countQuery1

$countQuery2 = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
        WHERE f.file_type = 'Aligned DNA reads file'
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@
$ws.Range("C3").Value = $countQuery2

$caseQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
WHERE f.file_type = 'Aligned DNA reads file' 
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@
$ws.Range("B2").Value = $caseQuery

$fileQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE f.file_type = 'Aligned DNA reads file'
 WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@
$ws.Range("B3").Value = $fileQuery

# ---- Remaining cells reuse pre-existing shared strings. ----
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

$ws.Range("D2").Value = "TC01_Trials_Filter_AssocFileType-AlignedDNA_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC01_Trials_Filter_AssocFileType-AlignedDNA_WebData.xlsx"

$ws.Range("D3").Value = "TC01_Trials_Filter_AssocFileType-AlignedDNA_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Trials_Filter_AssocFileType-AlignedDNA_WebData.xlsx"

# ---- Formatting: wrap-text style on B2, C2, B3, C3 (same style class as before) ----
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# ---- Row heights ----
$ws.Rows.Item(2).RowHeight = 188.5
$ws.Rows.Item(3).RowHeight = 409.5

# ---- Column widths (closest achievable to the target 75.8/70.3/28.5-ish
# ---- character widths given this engine's column-width quantisation) ----
$ws.Columns.Item(1).ColumnWidth = 8
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 27.666666667

# ---- View settings ----
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("C2").Select()
